$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are decimal-like strings (e.g. "1.001") that Excel would
# otherwise auto-coerce to numbers, corrupting formatting (e.g. "1.000" -> 1).
# Force the specific cells being rewritten to Text format first so the literal
# string is preserved, matching the source data feed which always stores these as text.
$priceCells = @('D2', 'D3', 'D4', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '23.220.25'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '1.602.27'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = '303.76'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('D7').Value = '0.3766'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = '52.01'
$ws.Range('E8').Value = '  +4.09%  '
$ws.Range('D9').Value = '0.3634'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = '1.273'
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.08146'
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '22.78'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '6.575'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Value = '7.410'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').Value = '0.00001250'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '1.601.10'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '94.08'
$ws.Range('E18').Value = '  +2.59%  '
$ws.Range('D19').Value = '0.06918'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').Value = '18.12'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').Value = '6.534'
$ws.Range('E21').Value = '  -0.36%  '
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = '12.90'
$ws.Range('E23').Value = '  -1.45%  '
$ws.Range('D24').Value = '23.212.78'
$ws.Range('D25').Value = '2.451'
$ws.Range('E25').Value = '  +3.77%  '
$ws.Range('D26').Value = '3.046'
$ws.Range('E26').Value = '  +7.09%  '
$ws.Range('D27').Value = '21.19'
$ws.Range('E27').Value = '  +0.73%  '
$ws.Range('D28').Value = '149.49'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').Value = '5.277'
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('D30').Value = '135.95'
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('D31').Value = '2.376'
$ws.Range('E31').Value = '  +5.85%  '
$ws.Range('D32').Value = '6.727'
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('D33').Value = '1.777.89'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').Value = '0.9646'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').Value = '0.07478'
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('D36').Value = '10.33'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').Value = '0.02740'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').Value = '0.2526'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('D39').Value = '6.126'
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').Value = '0.08784'
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D41').Value = '1.384'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').Value = '0.7097'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D43').Value = '12.42'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').Value = '15.64'
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('D45').Value = '0.6534'
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.315'
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').Value = '4.013'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').Value = '132.50'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = '0.07922'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').Value = '1.202'
$ws.Range('E51').Value = '  -2.82%  '
